# Apply "Natmi following Dr Hou advice" update to LR-pairs_lrc2p sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = New-Object 'object[,]' 15,20

# Row 2
$data[0,0] = "ECs"
$data[0,1] = "Apod"
$data[0,2] = "Lepr"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 6.320198333333333
$data[0,7] = 18.960595
$data[0,8] = 0.01706930024836835
$data[0,9] = 0.01706930024836834
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 7.613593333333334
$data[0,13] = 22.84078
$data[0,14] = 0.9151530262704457
$data[0,15] = 0.9151530262704456
$data[0,16] = 48.11941989601111
$data[0,17] = 433.0747790641
$data[0,18] = 0.01562102177861316
$data[0,19] = 0.01562102177861316

# Row 3
$data[1,0] = "ECs"
$data[1,1] = "Apod"
$data[1,2] = "Lepr"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 6.320198333333333
$data[1,7] = 18.960595
$data[1,8] = 0.01706930024836835
$data[1,9] = 0.01706930024836834
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.5558253333333333
$data[1,13] = 1.667476
$data[1,14] = 0.06681013991787221
$data[1,15] = 0.06681013991787221
$data[1,16] = 3.512926345357777
$data[1,17] = 31.61633710821999
$data[1,18] = 0.00114040233789366
$data[1,19] = 0.00114040233789366

# Row 4
$data[2,0] = "ECs"
$data[2,1] = "Apod"
$data[2,2] = "Lepr"
$data[2,3] = "sCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 6.320198333333333
$data[2,7] = 18.960595
$data[2,8] = 0.01706930024836835
$data[2,9] = 0.01706930024836834
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.150057
$data[2,13] = 0.450171
$data[2,14] = 0.01803683381168212
$data[2,15] = 0.01803683381168212
$data[2,16] = 0.948390001305
$data[2,17] = 8.535510011744998
$data[2,18] = 0.0003078761318615242
$data[2,19] = 0.0003078761318615241

# Row 5
$data[3,0] = "FAPs"
$data[3,1] = "Apod"
$data[3,2] = "Lepr"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 357.712545
$data[3,7] = 1073.137635
$data[3,8] = 0.966093548205577
$data[3,9] = 0.9660935482055769
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 7.613593333333334
$data[3,13] = 22.84078
$data[3,14] = 0.9151530262704457
$data[3,15] = 0.9151530262704456
$data[3,16] = 2723.4778478617
$data[3,17] = 24511.3006307553
$data[3,18] = 0.8841234343006865
$data[3,19] = 0.8841234343006863

# Row 6
$data[4,0] = "FAPs"
$data[4,1] = "Apod"
$data[4,2] = "Lepr"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 357.712545
$data[4,7] = 1073.137635
$data[4,8] = 0.966093548205577
$data[4,9] = 0.9660935482055769
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.5558253333333333
$data[4,13] = 1.667476
$data[4,14] = 0.06681013991787221
$data[4,15] = 0.06681013991787221
$data[4,16] = 198.82569456214
$data[4,17] = 1789.43125105926
$data[4,18] = 0.06454484512936823
$data[4,19] = 0.06454484512936821

# Row 7
$data[5,0] = "FAPs"
$data[5,1] = "Apod"
$data[5,2] = "Lepr"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 357.712545
$data[5,7] = 1073.137635
$data[5,8] = 0.966093548205577
$data[5,9] = 0.9660935482055769
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.150057
$data[5,13] = 0.450171
$data[5,14] = 0.01803683381168212
$data[5,15] = 0.01803683381168212
$data[5,16] = 53.677271365065
$data[5,17] = 483.095442285585
$data[5,18] = 0.0174252687755223
$data[5,19] = 0.0174252687755223

# Row 8
$data[6,0] = "M1"
$data[6,1] = "Apod"
$data[6,2] = "Lepr"
$data[6,3] = "ECs"
$data[6,4] = 1
$data[6,5] = 0.3333333333333333
$data[6,6] = 0.2484933333333333
$data[6,7] = 0.74548
$data[6,8] = 0.0006711193371913507
$data[6,9] = 0.0006711193371913505
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 7.613593333333334
$data[6,13] = 22.84078
$data[6,14] = 0.9151530262704457
$data[6,15] = 0.9151530262704456
$data[6,16] = 1.891927186044445
$data[6,17] = 17.0273446744
$data[6,18] = 0.0006141768924192802
$data[6,19] = 0.0006141768924192801

# Row 9
$data[7,0] = "M1"
$data[7,1] = "Apod"
$data[7,2] = "Lepr"
$data[7,3] = "FAPs"
$data[7,4] = 1
$data[7,5] = 0.3333333333333333
$data[7,6] = 0.2484933333333333
$data[7,7] = 0.74548
$data[7,8] = 0.0006711193371913507
$data[7,9] = 0.0006711193371913505
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.5558253333333333
$data[7,13] = 1.667476
$data[7,14] = 0.06681013991787221
$data[7,15] = 0.06681013991787221
$data[7,16] = 0.1381188898311111
$data[7,17] = 1.24307000848
$data[7,18] = 0.0000448375768193437994856009876087910015
$data[7,19] = 0.000044837576819343792709337409574388289

# Row 10
$data[8,0] = "M1"
$data[8,1] = "Apod"
$data[8,2] = "Lepr"
$data[8,3] = "sCs"
$data[8,4] = 1
$data[8,5] = 0.3333333333333333
$data[8,6] = 0.2484933333333333
$data[8,7] = 0.74548
$data[8,8] = 0.0006711193371913507
$data[8,9] = 0.0006711193371913505
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.150057
$data[8,13] = 0.450171
$data[8,14] = 0.01803683381168212
$data[8,15] = 0.01803683381168212
$data[8,16] = 0.03728816412
$data[8,17] = 0.33559347708
$data[8,18] = 0.0000121048679527266497997495253713395869
$data[8,19] = 0.000012104867952726639635354158319735518

# Row 11
$data[9,0] = "M2"
$data[9,1] = "Apod"
$data[9,2] = "Lepr"
$data[9,3] = "ECs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.2365366666666667
$data[9,7] = 0.7096100000000001
$data[9,8] = 0.0006388273231533433
$data[9,9] = 0.0006388273231533432
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 7.613593333333334
$data[9,13] = 22.84078
$data[9,14] = 0.9151530262704457
$data[9,15] = 0.9151530262704456
$data[9,16] = 1.800893988422223
$data[9,17] = 16.20804589580001
$data[9,18] = 0.0005846247580480301
$data[9,19] = 0.0005846247580480299

# Row 12
$data[10,0] = "M2"
$data[10,1] = "Apod"
$data[10,2] = "Lepr"
$data[10,3] = "FAPs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 0.2365366666666667
$data[10,7] = 0.7096100000000001
$data[10,8] = 0.0006388273231533433
$data[10,9] = 0.0006388273231533432
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 0.5558253333333333
$data[10,13] = 1.667476
$data[10,14] = 0.06681013991787221
$data[10,15] = 0.06681013991787221
$data[10,16] = 0.1314730715955556
$data[10,17] = 1.18325764436
$data[10,18] = 0.0000426801428432346276458514422458279114
$data[10,19] = 0.0000426801428432346276458514422458279114

# Row 13
$data[11,0] = "M2"
$data[11,1] = "Apod"
$data[11,2] = "Lepr"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 0.2365366666666667
$data[11,7] = 0.7096100000000001
$data[11,8] = 0.0006388273231533433
$data[11,9] = 0.0006388273231533432
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 0.150057
$data[11,13] = 0.450171
$data[11,14] = 0.01803683381168212
$data[11,15] = 0.01803683381168212
$data[11,16] = 0.03549398259000001
$data[11,17] = 0.31944584331
$data[11,18] = 0.0000115224222620785992249889742589097352
$data[11,19] = 0.0000115224222620785992249889742589097352

# Row 14
$data[12,0] = "sCs"
$data[12,1] = "Apod"
$data[12,2] = "Lepr"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 5.749211333333332
$data[12,7] = 17.247634
$data[12,8] = 0.01552720488570988
$data[12,9] = 0.01552720488570988
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 7.613593333333334
$data[12,13] = 22.84078
$data[12,14] = 0.9151530262704457
$data[12,15] = 0.9151530262704456
$data[12,16] = 43.77215707939111
$data[12,17] = 393.94941371452
$data[12,18] = 0.01420976854067865
$data[12,19] = 0.01420976854067864

# Row 15
$data[13,0] = "sCs"
$data[13,1] = "Apod"
$data[13,2] = "Lepr"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 5.749211333333332
$data[13,7] = 17.247634
$data[13,8] = 0.01552720488570988
$data[13,9] = 0.01552720488570988
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 0.5558253333333333
$data[13,13] = 1.667476
$data[13,14] = 0.06681013991787221
$data[13,15] = 0.06681013991787221
$data[13,16] = 3.195557305753777
$data[13,17] = 28.760015751784
$data[13,18] = 0.001037374730947746
$data[13,19] = 0.001037374730947746

# Row 16
$data[14,0] = "sCs"
$data[14,1] = "Apod"
$data[14,2] = "Lepr"
$data[14,3] = "sCs"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 5.749211333333332
$data[14,7] = 17.247634
$data[14,8] = 0.01552720488570988
$data[14,9] = 0.01552720488570988
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 0.150057
$data[14,13] = 0.450171
$data[14,14] = 0.01803683381168212
$data[14,15] = 0.01803683381168212
$data[14,16] = 0.8627094050459998
$data[14,17] = 7.764384645413998
$data[14,18] = 0.0002800616140834877
$data[14,19] = 0.0002800616140834877

$ws.Range("A2:T16").Value = $data

Write-Output "Applied rows 2-16 per Dr Hou revision"